$d = $word.ActiveDocument

$d.Content.Find.Execute("79×42=3318", $true, $false, $false, $false, $false, $true, 1, $false, "27×18=486", 2)
$d.Content.Find.Execute("81×26=2106", $true, $false, $false, $false, $false, $true, 1, $false, "77×36=2772", 2)
$d.Content.Find.Execute("50×37=1850", $true, $false, $false, $false, $false, $true, 1, $false, "32×47=1504", 2)
$d.Content.Find.Execute("83×56=4648", $true, $false, $false, $false, $false, $true, 1, $false, "23×90=2070", 2)
$d.Content.Find.Execute("73×68=4964", $true, $false, $false, $false, $false, $true, 1, $false, "82×57=4674", 2)
$d.Content.Find.Execute("73×58=4234", $true, $false, $false, $false, $false, $true, 1, $false, "56×79=4424", 2)
$d.Content.Find.Execute("26×43=1118", $true, $false, $false, $false, $false, $true, 1, $false, "66×51=3366", 2)
$d.Content.Find.Execute("29×22=638", $true, $false, $false, $false, $false, $true, 1, $false, "44×41=1804", 2)
$d.Content.Find.Execute("59×45=2655", $true, $false, $false, $false, $false, $true, 1, $false, "31×21=651", 2)
$d.Content.Find.Execute("40×38=1520", $true, $false, $false, $false, $false, $true, 1, $false, "60×14=840", 2)
$d.Content.Find.Execute("87×32=2784", $true, $false, $false, $false, $false, $true, 1, $false, "29×78=2262", 2)
$d.Content.Find.Execute("66×91=6006", $true, $false, $false, $false, $false, $true, 1, $false, "86×84=7224", 2)
$d.Content.Find.Execute("29×42=1218", $true, $false, $false, $false, $false, $true, 1, $false, "20×27=540", 2)
$d.Content.Find.Execute("13×30=390", $true, $false, $false, $false, $false, $true, 1, $false, "34×54=1836", 2)
$d.Content.Find.Execute("21×57=1197", $true, $false, $false, $false, $false, $true, 1, $false, "67×47=3149", 2)
$d.Content.Find.Execute("89×27=2403", $true, $false, $false, $false, $false, $true, 1, $false, "97×35=3395", 2)
$d.Content.Find.Execute("95×38=3610", $true, $false, $false, $false, $false, $true, 1, $false, "65×77=5005", 2)
$d.Content.Find.Execute("85×32=2720", $true, $false, $false, $false, $false, $true, 1, $false, "42×89=3738", 2)
$d.Content.Find.Execute("21×68=1428", $true, $false, $false, $false, $false, $true, 1, $false, "31×15=465", 2)
$d.Content.Find.Execute("59×99=5841", $true, $false, $false, $false, $false, $true, 1, $false, "13×15=195", 2)
$d.Content.Find.Execute("24×67=1608", $true, $false, $false, $false, $false, $true, 1, $false, "88×23=2024", 2)
$d.Content.Find.Execute("31×82=2542", $true, $false, $false, $false, $false, $true, 1, $false, "33×43=1419", 2)
$d.Content.Find.Execute("15×61=915", $true, $false, $false, $false, $false, $true, 1, $false, "70×69=4830", 2)
$d.Content.Find.Execute("13×69=897", $true, $false, $false, $false, $false, $true, 1, $false, "20×24=480", 2)
$d.Content.Find.Execute("25×31=775", $true, $false, $false, $false, $false, $true, 1, $false, "15×96=1440", 2)
